# Update raw flux data for siPHD2 experiments on the "dna" sheet,
# refresh sheet-tab selection / range selection to match the authored
# workbook state, and tag the active sheet/theme metadata.

$wb = $excel.ActiveWorkbook
$wsEvap = $wb.Worksheets.Item("evap")
$wsDna  = $wb.Worksheets.Item("dna")

# --- Updated raw counts (columns D/E/F, rows 2-57) on the "dna" sheet ----
$rawData = @"
2|44844|39883|41365
3|1058773|1058660|1014283
4|1942830|1929469|1942436
5|3927198|3791851|3790873
6|6810676|7392990|7307346
7|14981559|14346887|14270570
8|30379368|28899010|29143216
9|64868556|59817952|62458144
10|2133066|2188183|1916836
11|1091172|1771933|1364459
12|1106769|1329995|1441278
13|2600588|2505965|2039369
14|1648064|1669353|2459277
15|2494961|1984960|2646266
16|1547056|1266249|1148969
17|955681|774627|1232056
18|1621711|1048805|1198937
19|1268150|1562443|1517336
20|1509414|1447800|1520465
21|1140904|1623929|1388622
22|3454268|3815123|3646638
23|4335980|3867928|3695343
24|4337564|3642288|3647120
25|3414844|3771788|3961920
26|3392548|3585619|3807008
27|4072605|3943070|3869258
28|2149766|2099665|2079078
29|2387648|1510600|1855310
30|1717265|2143006|2793991
31|2036174|2107250|2682724
32|2044288|2254020|2250797
33|2013578|2224121|2304648
34|6352604|5799052|5876998
35|5467045|5401330|6047687
36|5573839|5969860|6258988
37|6074824|5657978|6895222
38|6709402|6107761|5974299
39|7127276|5717771|5819718
40|3835153|3385333|3891562
41|3946717|3153668|3887459
42|4584600|3819513|3795733
43|4247594|4264192|3667922
44|3214305|3753536|3610741
45|3227108|3512434|3796065
46|7365754|7615056|7610168
47|7675102|7408166|6502960
48|6394916|6087944|6723658
49|5861176|5595772|5255994
50|8071938|7979796|9742366
51|9368005|8806751|8545450
52|5084557|5153976|5448598
53|4227734|5167980|4467764
54|4587386|4173160|4726405
55|5440007|4662962|4261628
56|4522908|4512170|4653296
57|4543330|4508427|4267226
"@

$rows = $rawData -split "`n" | Where-Object { $_.Trim() -ne "" }
foreach ($line in $rows) {
    $parts = $line.Trim() -split '\|'
    $r = [int]$parts[0]
    $dVal = [double]$parts[1]
    $eVal = [double]$parts[2]
    $fVal = [double]$parts[3]

    $wsDna.Cells.Item($r, 4).Value = $dVal   # column D
    $wsDna.Cells.Item($r, 5).Value = $eVal   # column E
    $wsDna.Cells.Item($r, 6).Value = $fVal   # column F
}

# --- View / selection state -------------------------------------------
# "dna" was the active tab before; "evap" becomes the active tab, and the
# selection on "dna" moves down to the newly added block below the table.
$wsDna.Activate()
$wsDna.Range("D58:F65").Select()
$wsEvap.Activate()
$wsEvap.Range("C12").Select()

# --- Workbook theme rename (cosmetic, best-effort) ----------------------
try {
    $wb.Theme.Name = "Office Theme 2013 - 2022"
} catch {
}

Write-Output "applied siPHD2 raw data update"
